$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.967.91'
$ws.Range('E2').Value = '  +0.06%  '
$ws.Range('D3').Value = '2.402.60'
$ws.Range('E3').Value = '  -0.41%  '
$ws.Range('D5').Value = "'566.67"
$ws.Range('E5').Value = '  -0.22%  '
$ws.Range('D6').Value = "'142.02"
$ws.Range('E6').Value = '  +2.58%  '
$ws.Range('D7').Value = "'0.999"
$ws.Range('E7').Value = '  -0.41%  '
$ws.Range('D8').Value = "'0.540"
$ws.Range('E8').Value = '  +2.82%  '
$ws.Range('D9').Value = '2.410.09'
$ws.Range('E9').Value = '  +0.66%  '
$ws.Range('E10').Value = '  +2.21%  '
$ws.Range('E11').Value = '  +0.03%  '
$ws.Range('D12').Value = "'5.20"
$ws.Range('E12').Value = '  +2.97%  '
$ws.Range('E13').Value = '  +2.91%  '
$ws.Range('D14').Value = "'26.50"
$ws.Range('E15').Value = '  +0.37%  '
$ws.Range('D16').Value = '2.838.75'
$ws.Range('E16').Value = '  -0.63%  '
$ws.Range('D17').Value = '60.772.17'
$ws.Range('E17').Value = '  +0.05%  '
$ws.Range('D18').Value = '2.412.79'
$ws.Range('E18').Value = '  +0.22%  '
$ws.Range('E19').Value = '  +3.42%  '
$ws.Range('D20').Value = "'10.71"
$ws.Range('E20').Value = '  +1.21%  '
$ws.Range('D21').Value = "'324.51"
$ws.Range('E21').Value = '  +0.81%  '
$ws.Range('E22').Value = '  +1.70%  '
$ws.Range('D23').Value = "'6.04"
$ws.Range('E23').Value = '  -0.70%  '
$ws.Range('E24').Value = '  -0.24%  '
$ws.Range('D25').Value = "'1.90"
$ws.Range('E25').Value = '  +5.47%  '
$ws.Range('D26').Value = "'65.16"
$ws.Range('E26').Value = '  +0.62%  '
$ws.Range('D27').Value = "'585.33"
$ws.Range('E27').Value = '  +1.47%  '
$ws.Range('D28').Value = "'8.26"
$ws.Range('E28').Value = '  +1.22%  '
$ws.Range('B29').Value = 'PEPE'
$ws.Range('C29').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D29').Value = '0.0₃0947'
$ws.Range('E29').Value = '  +3.54%  '
$ws.Range('B30').Value = 'WrappedeETH'
$ws.Range('C30').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D30').Value = '2.518.47'
$ws.Range('E30').Value = '  -0.70%  '
$ws.Range('D31').Value = "'8.05"
$ws.Range('E31').Value = '  +3.16%  '
$ws.Range('E32').Value = '  +1.74%  '
$ws.Range('E33').Value = '  +0.08%  '
$ws.Range('D34').Value = "'0.133"
$ws.Range('E34').Value = '  +0.95%  '
$ws.Range('E35').Value = '  +5.91%  '
$ws.Range('E36').Value = '  -0.62%  '
$ws.Range('D37').Value = "'153.94"
$ws.Range('E37').Value = '  +1.22%  '
$ws.Range('E38').Value = '  +1.94%  '
$ws.Range('E39').Value = '  +1.54%  '
$ws.Range('D40').Value = "'18.33"
$ws.Range('E40').Value = '  +0.91%  '
$ws.Range('D41').Value = "'5.21"
$ws.Range('E42').Value = '  -0.13%  '
$ws.Range('E43').Value = '  +1.94%  '
$ws.Range('E44').Value = '  +11.89%  '
$ws.Range('D45').Value = "'41.67"
$ws.Range('E45').Value = '  +1.44%  '
$ws.Range('D46').Value = '0.0₆0286'
$ws.Range('E46').Value = '  +9.04%  '
$ws.Range('D47').Value = "'141.80"
$ws.Range('E47').Value = '  -0.07%  '
$ws.Range('D48').Value = "'3.53"
$ws.Range('E48').Value = '  +1.24%  '
$ws.Range('E49').Value = '  +1.55%  '
$ws.Range('D50').Value = "'0.0511"
$ws.Range('E50').Value = '  +2.31%  '
$ws.Range('D51').Value = "'19.57"
$ws.Range('E51').Value = '  +1.95%  '
